$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename result "Inconclusive" to "Non-Negative"
$ws.Range("B10").Value = "Non-Negative"

# Update the active cell selection to B11 (as saved by the author)
$ws.Range("B11").Select()
